# Update crypto price (Price) and 1h volume change (Volume(1h)) columns
# to reflect the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.775.61'
$ws.Range('E2').Value = '  -0.14%  '
$ws.Range('D3').Value = '2.547.91'
$ws.Range('E3').Value = '  +0.44%  '
$ws.Range('D4').Value = '''0.998'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '''319.28'
$ws.Range('E5').Value = '  +4.78%  '
$ws.Range('D6').Value = '''95.48'
$ws.Range('E6').Value = '  -2.52%  '
$ws.Range('D7').Value = '''0.579'
$ws.Range('E7').Value = '  +0.41%  '
$ws.Range('D9').Value = '''0.535'
$ws.Range('E9').Value = '  -1.95%  '
$ws.Range('D10').Value = '''36.39'
$ws.Range('E10').Value = '  -0.15%  '
$ws.Range('E11').Value = '  -1.22%  '
$ws.Range('D12').Value = '''7.73'
$ws.Range('E12').Value = '  +1.85%  '
$ws.Range('E13').Value = '  -0.53%  '
$ws.Range('D14').Value = '2.938.84'
$ws.Range('E14').Value = '  +0.52%  '
$ws.Range('D15').Value = '''15.95'
$ws.Range('E15').Value = '  +5.83%  '
$ws.Range('D16').Value = '2.538.60'
$ws.Range('E16').Value = '  +0.48%  '
$ws.Range('D17').Value = '''0.866'
$ws.Range('E17').Value = '  -0.46%  '
$ws.Range('D18').Value = '42.800.37'
$ws.Range('E18').Value = '  -0.10%  '
$ws.Range('D19').Value = '''13.14'
$ws.Range('E19').Value = '  -0.42%  '
$ws.Range('E20').Value = '  +1.33%  '
$ws.Range('E21').Value = '  -1.91%  '
$ws.Range('D22').Value = '''71.21'
$ws.Range('E22').Value = '  -0.54%  '
$ws.Range('D23').Value = '''253.41'
$ws.Range('E23').Value = '  -0.17%  '
$ws.Range('D24').Value = '''2.99'
$ws.Range('E24').Value = '  +2.32%  '
$ws.Range('D25').Value = '''2.03'
$ws.Range('E25').Value = '  -2.13%  '
$ws.Range('D26').Value = '''27.24'
$ws.Range('E26').Value = '  -1.79%  '
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('D28').Value = '''2.36'
$ws.Range('E28').Value = '  +3.32%  '
$ws.Range('D29').Value = '''40.05'
$ws.Range('E29').Value = '  +4.85%  '
$ws.Range('E30').Value = '  +0.99%  '
$ws.Range('D31').Value = '''6.00'
$ws.Range('E31').Value = '  -3.39%  '
$ws.Range('D32').Value = '''156.10'
$ws.Range('E32').Value = '  -0.70%  '
$ws.Range('E33').Value = '  +1.16%  '
$ws.Range('E34').Value = '  +1.90%  '
$ws.Range('D35').Value = '''19.23'
$ws.Range('E35').Value = '  -1.05%  '
$ws.Range('D36').Value = '''0.0794'
$ws.Range('E36').Value = '  -0.32%  '
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('E38').Value = '  -2.66%  '
$ws.Range('D39').Value = '''2.42'
$ws.Range('E39').Value = '  +13.76%  '
$ws.Range('E40').Value = '  -0.15%  '
$ws.Range('D41').Value = '''23.84'
$ws.Range('E41').Value = '  -4.12%  '
$ws.Range('D42').Value = '''3.87'
$ws.Range('E42').Value = '  +0.54%  '
$ws.Range('D43').Value = '''3.38'
$ws.Range('E43').Value = '  -1.21%  '
$ws.Range('E44').Value = '  +0.41%  '
$ws.Range('E45').Value = '  -0.68%  '
$ws.Range('D46').Value = '2.033.87'
$ws.Range('E46').Value = '  -3.14%  '
$ws.Range('D47').Value = '''84.86'
$ws.Range('E47').Value = '  -2.11%  '
$ws.Range('D48').Value = '''8.95'
$ws.Range('E48').Value = '  -0.20%  '
$ws.Range('D49').Value = '2.792.30'
$ws.Range('E49').Value = '  +0.33%  '
$ws.Range('D50').Value = '''74.19'
$ws.Range('E50').Value = '  +0.77%  '
$ws.Range('E51').Value = '  -0.42%  '
